# Update the akkupack-ng BoM workbook:
#  - Add 4 more diode references (D203 D204 D205 D206) to the "D" group row
#  - Add 4 more resistor references (R214 R215 R216 R217) to the "R" group row
#    and bump its "Quantity Per PCB" count from 38 to 42
#  - Update the component-count summary figures on both the BoM and DNF sheets

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- BoM sheet: group rows -------------------------------------------------
$bom.Range("B14").Value = "D102 D103 D104 D105 D106 D107 D108 D109 D110 D111 D203 D204 D205 D206 D301 D401"

$bom.Range("B24").Value = "R101 R102 R103 R104 R105 R106 R107 R108 R109 R110 R111 R112 R113 R114 R201 R202 R203 R204 R205 R206 R207 R208 R209 R210 R211 R212 R213 R214 R215 R216 R217 R301 R302 R401 R402 R501 R502 R503 R504 R505 R506 R507"
$bom.Range("L24").Value = 42

# --- Summary counters on BoM sheet -----------------------------------------
$bom.Range("F3").Value = 127
$bom.Range("F4").Value = 126
$bom.Range("F6").Value = 126

# --- Summary counters on DNF sheet ------------------------------------------
$dnf.Range("F3").Value = 127
$dnf.Range("F4").Value = 126
$dnf.Range("F6").Value = 126
